# Error Calculations and Plots
# Apply the missing-data shuffle: two whole rows ("RM 232" and "SC 92") are
# dropped from the sheet (rows shift up), and a number of individual cells
# flip between a value and blank across the remaining rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the two rows that disappear entirely (delete the lower row first
#    so the earlier row's index stays valid).
$ws.Rows(28).Delete()   # "SC 92"
$ws.Rows(26).Delete()   # "RM 232"

# 2) Apply the remaining per-cell value changes (post row-shift row numbers).
$ws.Cells.Item(2, 3).Value = 14.9

$ws.Cells.Item(3, 3).Value = ""

$ws.Cells.Item(4, 3).Value = ""

$ws.Cells.Item(5, 5).Value = ""

$ws.Cells.Item(6, 6).Value = 16.43

$ws.Cells.Item(8, 5).Value = -6.6

$ws.Cells.Item(10, 5).Value = -6.1

$ws.Cells.Item(11, 3).Value = 11.4
$ws.Cells.Item(11, 6).Value = 17.65

$ws.Cells.Item(12, 5).Value = ""
$ws.Cells.Item(12, 6).Value = ""

$ws.Cells.Item(13, 3).Value = ""
$ws.Cells.Item(13, 6).Value = 17.1

$ws.Cells.Item(15, 5).Value = -8.4

$ws.Cells.Item(17, 6).Value = ""

$ws.Cells.Item(18, 5).Value = ""
$ws.Cells.Item(18, 6).Value = 18.35

$ws.Cells.Item(19, 5).Value = ""
$ws.Cells.Item(19, 6).Value = ""

$ws.Cells.Item(21, 3).Value = 12.7

$ws.Cells.Item(24, 6).Value = ""

$ws.Cells.Item(25, 3).Value = ""
$ws.Cells.Item(25, 5).Value = -7.1
$ws.Cells.Item(25, 6).Value = 16.6

$ws.Cells.Item(27, 5).Value = -10

$ws.Cells.Item(29, 2).Value = ""
$ws.Cells.Item(29, 5).Value = ""

$ws.Cells.Item(31, 6).Value = ""

$ws.Cells.Item(32, 6).Value = ""

$ws.Cells.Item(33, 2).Value = -19.5
$ws.Cells.Item(33, 3).Value = 10.4
$ws.Cells.Item(33, 5).Value = ""

Write-Host "Edit applied"
